$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Homework 2" header in G2
$ws.Range("G2").Value = "Homework 2"

# Set formulas for the new Homework 2 column (out of 12 points)
$ws.Range("G3").Formula = "=10.5/12"
$ws.Range("G4").Formula = "=6/12"
$ws.Range("G5").Formula = "=12/12"
$ws.Range("G6").Formula = "=12/12"
$ws.Range("G7").Formula = "=12/12"
$ws.Range("G8").Formula = "=12/12"
$ws.Range("G9").Formula = "=10/12"
$ws.Range("G11").Formula = "=12/12"
$ws.Range("G12").Formula = "=11/12"
$ws.Range("G13").Formula = "=12/12"
$ws.Range("G14").Formula = "=12/12"
$ws.Range("G15").Formula = "=12/12"
$ws.Range("G16").Formula = "=12/12"
$ws.Range("G17").Formula = "=12/12"
$ws.Range("G18").Formula = "=12/12"
$ws.Range("G19").Formula = "=12/12"

# Adjust column widths: F widened, G added with similar width to B
$ws.Columns.Item(6).ColumnWidth = 12.91
$ws.Columns.Item(7).ColumnWidth = 13.11

# Update the active cell selection to G3
$ws.Range("G3").Select()
